$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Reuse the exact formatting (number format / alignment) already present on
# row 2 for the four new rows by copying its formats down, then fill in the
# new values. This keeps the same style indices (s="1","2","3","4") instead
# of minting duplicate styleSheet entries. Row 6 has no Notes entry, so only
# stamp the Notes-column (D) format on rows 3-5.
$ws.Range("A2:C2").Copy()
$ws.Range("A3:C6").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D3:D5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 3 - "Notebook" entry
$ws.Cells.Item(3, 1).Value = 41759
$ws.Cells.Item(3, 2).Value = 0.44861111111111113
$ws.Cells.Item(3, 3).Value = "Notebook"
$ws.Cells.Item(3, 4).Value = "Started this notebook. Created git repo and started syncing. One issue with excel is that it locks and must close to commit to repo."
$ws.Rows.Item(3).RowHeight = 30

# Row 4 - "DOL - Venture" entry
$ws.Cells.Item(4, 1).Value = 41759
$ws.Cells.Item(4, 2).Value = 0.45208333333333334
$ws.Cells.Item(4, 3).Value = "DOL - Venture"
$ws.Cells.Item(4, 4).Value = "Steps to locate sproc. To be continued."

# Row 5 - "Algorithms I" entry
$ws.Cells.Item(5, 1).Value = 41759
$ws.Cells.Item(5, 2).Value = 0.45555555555555555
$ws.Cells.Item(5, 3).Value = "Algorithms I"
$ws.Cells.Item(5, 4).Value = "Continued from last time.  Assumption made is that f(n) and g(n) are always positive numbers. That condition should find big O. Then we need big omega… between those two should be theta. So, for omega: 2(max(f(n), g(n))) >= f(n) + g(n) --> max(f(n), g(n)) >=  0.5(f(n) + g(n)). Thus: 1/2(f(n) + g(n)) <= max(f(n), g(n)) <= f(n) + g(n), for every n.  So, to satisfy theta, our constant C must be between 1 and 1/2... and n sub 0 is 1.`nNew vid: O(n log n) Algorithm for Counting Inversions.`n1. Divide prob into smaller subproblems`n2. Subproblems solved via recursion`n3. Combine subproblem solutions into one solution for the real problem.`nProblem: Input of array A containing numbers in some arbitrary order. Output is the number of inversions (number of pairs where indices i and j is such that i < j and A[i] > A[j]. (sorted order has 0 inversions, but any other array with have some non-zero number"
$ws.Rows.Item(5).RowHeight = 195

# Row 6 - "DOL - Venture" entry (no notes this time)
$ws.Cells.Item(6, 1).Value = 41759
$ws.Cells.Item(6, 2).Value = 0.49444444444444446
$ws.Cells.Item(6, 3).Value = "DOL - Venture"

$ws.Range("A7").Select()
